$wb = $excel.ActiveWorkbook

# Sheet index 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1225  # was 1222
$ws.Range("F11").Value = 2220  # was 2217
$ws.Range("F12").Value = 1552  # was 1551
$ws.Range("F13").Value = 1238  # was 1235
$ws.Range("F15").Value = 222  # was 221
$ws.Range("F17").Value = 714  # was 713
$ws.Range("F18").Value = 265  # was 264
$ws.Range("F19").Value = 1080  # was 1079
$ws.Range("F22").Value = 4182  # was 4176
$ws.Range("F24").Value = 143  # was 142
$ws.Range("F25").Value = 118  # was 117
$ws.Range("F28").Value = 600  # was 601
$ws.Range("F34").Value = 919  # was 918
$ws.Range("F37").Value = 112  # was 110

# Sheet index 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 775  # was 773
$ws.Range("F5").Value = 418  # was 417

# Sheet index 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1225  # was 1222
$ws.Range("F4").Value = 775  # was 773
$ws.Range("F8").Value = 418  # was 417
$ws.Range("F16").Value = 2220  # was 2217
$ws.Range("F17").Value = 1552  # was 1551
$ws.Range("F18").Value = 1238  # was 1235
$ws.Range("F20").Value = 222  # was 221
$ws.Range("F23").Value = 714  # was 713
$ws.Range("F24").Value = 265  # was 264
$ws.Range("F25").Value = 1080  # was 1079
$ws.Range("F28").Value = 4182  # was 4176
$ws.Range("F30").Value = 143  # was 142
$ws.Range("F31").Value = 118  # was 117
$ws.Range("F34").Value = 600  # was 601
$ws.Range("F40").Value = 919  # was 918
$ws.Range("F43").Value = 112  # was 110
